$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:group-code" (column C) and "codeforiati:group-name"
# (column D) values need to swap places, for the header row and every
# data row, so that the group-name now comes before the group-code.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$srcRange = $ws.Range("C1:D$lastRow")
$vals = $srcRange.Value2

$rows = $vals.GetLength(0)
$cols = $vals.GetLength(1)

$swapped = New-Object 'object[,]' $rows,$cols
for ($r = 1; $r -le $rows; $r++) {
    $swapped[$r-1,0] = $vals[$r,2]
    $swapped[$r-1,1] = $vals[$r,1]
}

$srcRange.Value = $swapped
